# Insert a new data record as row 229 in the "Mango" price table, shifting
# every existing record from row 229 down by one (old row 336 becomes 337).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 229..336 down to 230..337, carrying formatting (incl. the date
# number-format on column D) from the row being pushed down.
$ws.Rows("229:229").Insert()

# Populate the newly opened row 229 with the new market record.
$ws.Range("A229").Value = 4
$ws.Range("B229").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C229").Value = "Los Lagos"
$ws.Range("D229").Value = 44992
$ws.Range("E229").Value = 10
$ws.Range("F229").Value = "Fruta"
$ws.Range("G229").Value = 100108
$ws.Range("H229").Value = "Tropicales y subtropicales"
$ws.Range("I229").Value = 100108002
$ws.Range("J229").Value = "Mango"
$ws.Range("K229").Value = "Sin especificar"
$ws.Range("L229").Value = "Primera"
$ws.Range("M229").Value = 200
$ws.Range("N229").Value = 8500
$ws.Range("O229").Value = 9000
$ws.Range("P229").Value = 8750
$ws.Range("Q229").Value = "`$/bandeja 4 kilos"
$ws.Range("R229").Value = "Perú"
$ws.Range("S229").Value = 2188
$ws.Range("T229").Value = 4
